$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.492.48"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.886.02"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "243.92"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4722"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "0.2888"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "0.06478"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").Value = "22.24"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("D11").Value = "0.07761"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "1.886.99"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "95.73"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "0.7263"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").Value = "5.187"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "281.88"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "30.500.55"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "13.07"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "0.000007462"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "2.137.40"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "5.272"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "6.308"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").Value = "163.88"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").Value = "9.070"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("D27").Value = "18.86"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").Value = "1.892"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").Value = "1.334"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "0.09652"
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("E31").Value = "  -2.35%  "
$ws.Range("D32").Value = "4.272"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").Value = "4.146"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").Value = "0.04852"
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").Value = "0.6928"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "0.01884"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").Value = "2.824"
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("D40").Value = "74.78"
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("D41").Value = "6.213"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").Value = "1.964"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "0.4267"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "0.8266"
$ws.Range("D46").Value = "101.09"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "9.639"
$ws.Range("E47").Value = "  +3.29%  "
$ws.Range("D48").Value = "6.959"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "35.19"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "906.67"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "0.05755"
$ws.Range("E51").Value = "  +1.74%  "
